# Add the "Reused Terms" sheet (Reuse Existing Terms metric) after "Short URI"
# and populate it with the metric's data, mirroring the structure used by the
# other metric sheets in this workbook.

$wb = $excel.ActiveWorkbook

$shortUri = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $shortUri)
$ws.Name = "Reused Terms"

# ---- Classes block (column A, filled in first) ------------------------------
$ws.Range("A2").Value = "Classes"
$ws.Range("A2").Font.Bold = $true

$ws.Range("A3").Value = "swrc:Conference"
$ws.Range("D3").Value = 1

$ws.Range("A4").Value = "geonames:SpatialThing"
$ws.Range("C4").Value = 0

# ---- Properties block (column A continues) ----------------------------------
$ws.Range("A6").Value = "Properties"
$ws.Range("A6").Font.Bold = $true

$ws.Range("A7").Value = "rdf:type"
$ws.Range("F7").Value = 1

$ws.Range("A8").Value = "swrc:startDate"
$ws.Range("D8").Value = 0

$ws.Range("A9").Value = "swrc:description"
$ws.Range("D9").Value = 0

$ws.Range("A10").Value = "rdfs:label"
$ws.Range("F10").Value = 1

$ws.Range("A11").Value = "dcterms:spatial"
$ws.Range("F11").Value = 1

$ws.Range("A12").Value = "swrc:eventTitle"
$ws.Range("D12").Value = 1

$ws.Range("A13").Value = "swrc:location"
$ws.Range("D13").Value = 1

$ws.Range("A14").Value = "owl:sameAs"
$ws.Range("F14").Value = 1

$ws.Range("A15").Value = "rdfs:seeAlso"
$ws.Range("F15").Value = 1

$ws.Range("A16").Value = "geonames:name"
$ws.Range("C16").Value = 1

$ws.Range("A17").Value = "geo:long"
$ws.Range("B17").Value = 1

$ws.Range("A18").Value = "geonames:P"
$ws.Range("C18").Value = 1

$ws.Range("A19").Value = "geo:lat"
$ws.Range("B19").Value = 1

$ws.Range("A20").Value = "geonames:countryName"
$ws.Range("C20").Value = 0

# ---- Column headers (row 1) : per-namespace totals, filled in afterwards ---
$ws.Range("B1").Value = "wgs84_pos"
$ws.Range("C1").Value = "geonames"
$ws.Range("D1").Value = "swrc"
$ws.Range("E1").Value = "swc"
$ws.Range("F1").Value = "Others (Top)"

# ---- Totals row --------------------------------------------------------------
$ws.Range("A22").Value = "Total"
$ws.Range("A22").Font.Bold = $true

$ws.Range("B22").Formula = "=COUNTIF(B2:B20,1)"
$ws.Range("C22").Formula = "=COUNTIF(C2:C20,1)"
$ws.Range("D22").Formula = "=COUNTIF(D2:D20,1)"
$ws.Range("E22").Formula = "=COUNTIF(E2:E20,1)"
$ws.Range("F22").Formula = "=COUNTIF(F2:F20,1)"

# ---- Summary box (H/I columns) ----------------------------------------------
$ws.Range("H6").Value = "Total Reused Terms"
$ws.Range("I6").Formula = "=SUM(B22:F22)"

$ws.Range("H7").Value = "Total Unique Terms"
$ws.Range("I7").Value = 16

$ws.Range("H8").Value = "Metric Value"
$ws.Range("I8").Formula = "=I6/I7"

# ---- Column widths, matching the other metric sheets -----------------------
$ws.Columns.Item(1).ColumnWidth = 27.998697916666668
$ws.Columns.Item(6).ColumnWidth = 10.666666666666666
$ws.Columns.Item(8).ColumnWidth = 16.498697916666668

# ---- Selection / active sheet -----------------------------------------------
$ws.Range("F12").Select()
$ws.Activate()
